$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "1+79=",
    "90-24=",
    "14+46=",
    "61-49=",
    "9-6=",
    "17+52=",
    "90-89=",
    "60-36=",
    "82-3=",
    "78-58=",
    "90+9=",
    "34+18=",
    "97-72=",
    "47-38=",
    "70-10=",
    "96-46=",
    "24+41=",
    "73-44=",
    "88-23=",
    "52+37=",
    "3+71=",
    "96-5=",
    "45+41=",
    "1+55=",
    "42-20=",
    "57-18=",
    "87-11=",
    "79+9=",
    "44+10=",
    "46-34=",
    "64-61=",
    "93-33=",
    "42-14=",
    "55+31=",
    "83-82=",
    "82-65=",
    "76+10=",
    "93-0=",
    "85-39=",
    "77-16=",
    "16+58=",
    "52+3=",
    "54-1=",
    "81-31=",
    "81-19=",
    "51+44=",
    "55+44=",
    "36+6=",
    "79-7=",
    "29-26=",
    "41-2=",
    "55+35=",
    "51-11=",
    "1+16=",
    "69-10=",
    "97-85=",
    "42+5=",
    "37+21=",
    "18+80=",
    "62-22=",
    "36+14=",
    "41+47=",
    "84+0=",
    "17+34=",
    "12+87=",
    "74-30=",
    "66-29=",
    "39+15=",
    "79+13=",
    "86-18=",
    "33-5=",
    "24-2=",
    "74+23=",
    "44+38=",
    "40-10=",
    "92-5=",
    "19+31=",
    "34-7=",
    "64-36=",
    "62+34=",
    "31+68=",
    "50-37=",
    "96-25=",
    "80+9=",
    "97-25=",
    "86-68=",
    "44-21=",
    "43-25=",
    "65+20=",
    "37+19=",
    "44-2=",
    "18-0=",
    "78-57=",
    "96-26=",
    "51+17=",
    "22+62=",
    "48-45=",
    "62-54=",
    "64-9=",
    "36-5="
)

$cols = 5
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated " + $idx + " cells")